$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: "noobject" ---
# Copy header formatting from F1 (same style as the other header cells), then set the text.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "noobject"

# Formula column: G2 standalone, G3:G13 filled as one batch (creates the shared formula group
# starting at G3, matching a manual fill-down from G2).
$ws.Range("G2").Formula = "=E2-F2"
$ws.Range("G3:G13").Formula = "=E3-F3"

# Apply the 2-decimal number format to the new formula cells.
$ws.Range("G2:G13").NumberFormat = "0.00"

# --- Clean up the old stray empty formatted cells scattered around H:Q ---
$ws.Range("I13:Q13").Clear()
$ws.Range("P14:Q14").Clear()
$ws.Range("Q15").Clear()
$ws.Range("H16:Q16").Clear()
$ws.Range("H17:I17").Clear()
$ws.Range("P17:Q17").Clear()
$ws.Range("H18:I18").Clear()
$ws.Range("P18:Q18").Clear()
$ws.Range("H19:Q19").Clear()
$ws.Range("H20:J20").Clear()
$ws.Range("Q20").Clear()
$ws.Range("H21").Clear()
$ws.Range("H22").Clear()

# Update the selection to match the filled formula range.
$ws.Range("G2:G13").Select() | Out-Null
